# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the file that was just handed back (6f3c7c99-3344-4975-854e-52c6f46f233c),
# and roll that up into the Overview sheet's "Latest HO Xliff Generate Date".

$wb = $excel.ActiveWorkbook

# zh-cn: new Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 12:45:36"
$wsZhCn.Range("K2").Value = "2016-08-27 12:45:53"

# de-de: new Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 12:45:40"
$wsDeDe.Range("K2").Value = "2016-08-27 12:45:59"

# Overview: Latest HO Xliff Generate Date is the max handoff datetime across languages
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 12:45:40"
